$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

$ws.Range("A2").Value = "yam"
$ws.Range("B2").Value = "Yambit12"
$ws.Range("C2").Value = "sjhfdasfhn"
$ws.Range("D2").Value = "Y@123456"

$ws.Range("A3").Value = "yam"
$ws.Range("B3").Value = "yambit33"
$ws.Range("C3").Value = "fafasf"
$ws.Range("D3").Value = "Y@123456"
